$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '318.72'
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '4.33%'
$ws.Cells.Item(2, 5).ClearFormats()
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '39.74'
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '2.40%'
$ws.Cells.Item(3, 5).ClearFormats()
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '0.55%'
$ws.Cells.Item(4, 5).ClearFormats()
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.08216'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '1.89%'
$ws.Cells.Item(5, 5).ClearFormats()
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '2.056'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '6.62%'
$ws.Cells.Item(6, 5).ClearFormats()
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '8.372'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '4.40%'
$ws.Cells.Item(7, 5).ClearFormats()
$ws.Cells.Item(8, 2).Value = 'MXToken'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.9425'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '1.69%'
$ws.Cells.Item(8, 5).ClearFormats()
$ws.Cells.Item(9, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.1368'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '-5.15%'
$ws.Cells.Item(9, 5).ClearFormats()
$ws.Cells.Item(10, 2).Value = 'WazirX'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.2002'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '4.74%'
$ws.Cells.Item(10, 5).ClearFormats()
$ws.Cells.Item(11, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.09097'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '0.71%'
$ws.Cells.Item(11, 5).ClearFormats()
$ws.Cells.Item(12, 2).Value = 'BitrueCoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.03525'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '0.44%'
$ws.Cells.Item(12, 5).ClearFormats()
$ws.Cells.Item(13, 2).Value = 'BitMartToken'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.09797'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '0.19%'
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(14, 2).Value = 'BitForexToken'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.001420'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '0.80%'
$ws.Cells.Item(14, 5).ClearFormats()
$ws.Cells.Item(15, 2).Value = 'TigerCash'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.006270'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '5.44%'
$ws.Cells.Item(15, 5).ClearFormats()
$ws.Cells.Item(16, 2).Value = 'LEO'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.683'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '-2.02%'
$ws.Cells.Item(16, 5).ClearFormats()
$ws.Cells.Item(17, 2).Value = 'GateToken'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '4.318'
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '2.73%'
$ws.Cells.Item(17, 5).ClearFormats()
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.358'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '0.87%'
$ws.Cells.Item(18, 5).ClearFormats()
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.3472'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '0.34%'
$ws.Cells.Item(19, 5).ClearFormats()
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.1321'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '2.21%'
$ws.Cells.Item(20, 5).ClearFormats()
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.955'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '5.69%'
$ws.Cells.Item(21, 5).ClearFormats()
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.2451'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '1.53%'
$ws.Cells.Item(22, 5).ClearFormats()
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.04345'
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '-0.59%'
$ws.Cells.Item(23, 5).ClearFormats()
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.001230'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '1.59%'
$ws.Cells.Item(24, 5).ClearFormats()
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.004793'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '12.40%'
$ws.Cells.Item(25, 5).ClearFormats()
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '-0.08%'
$ws.Cells.Item(26, 5).ClearFormats()
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '-10.05%'
$ws.Cells.Item(27, 5).ClearFormats()
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.02338'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '15.05%'
$ws.Cells.Item(39, 5).ClearFormats()
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.05203'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '2.93%'
$ws.Cells.Item(40, 5).ClearFormats()
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.007758'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '3.16%'
$ws.Cells.Item(41, 5).ClearFormats()
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.009915'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '1.26%'
$ws.Cells.Item(42, 5).ClearFormats()
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.1405'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '4.97%'
$ws.Cells.Item(43, 5).ClearFormats()
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.002080'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '-0.70%'
$ws.Cells.Item(44, 5).ClearFormats()
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.008922'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '-9.31%'
$ws.Cells.Item(45, 5).ClearFormats()
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.00006621'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '6.70%'
$ws.Cells.Item(46, 5).ClearFormats()
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '-0.20%'
$ws.Cells.Item(47, 5).ClearFormats()
$ws.Cells.Item(48, 2).Value = 'BOLO'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.002946'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '2.48%'
$ws.Cells.Item(48, 5).ClearFormats()
$ws.Cells.Item(49, 2).Value = 'CoinbaseStockToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.001690'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '-6.28%'
$ws.Cells.Item(49, 5).ClearFormats()
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.00002100'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '-0.20%'
$ws.Cells.Item(50, 5).ClearFormats()
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0002000'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '-0.20%'
$ws.Cells.Item(51, 5).ClearFormats()
